# Fix sheet name in test data `input.xlsx`
#
# The "features-feature1" sheet name used a dash where the i18n import
# logic expects a dot-separated "group.key" style name, so rename it to
# "features.feature-one".
$wb = $excel.ActiveWorkbook

$featuresSheet = $wb.Worksheets.Item("features-feature1")
$featuresSheet.Name = "features.feature-one"

# Bring the renamed sheet to the front / make it the active tab and set
# its selection, matching the interactive session in which the rename
# was made and verified.
$featuresSheet.Activate()
$featuresSheet.Range("C3").Select()
